$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PINS")

$ws.Range("B7").Value = 184000000.0
$ws.Range("C7").Value = 173000000.0
$ws.Range("D7").Value = 142757000.0
$ws.Range("E7").Value = 98580000.0
$ws.Range("F7").Value = 71914000.0
